$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Remove D128 (value dropped from this row in the new export) ---
$ws.Range("D128").ClearContents()

# --- Updated recomputed error values for existing rows ---
$ws.Range("D129").Value = 0.7804486234241067
$ws.Range("D130").Value = 0.7918746594241067
$ws.Range("D131").Value = 0.7115302104241067
$ws.Range("D132").Value = 0.7732527034241068
$ws.Range("C133").Value = 0.5947585844241068
$ws.Range("C134").Value = -0.1471494035758933
$ws.Range("C135").Value = 0.2315426864241067
$ws.Range("C136").Value = 0.5186180304241067
$ws.Range("C137").Value = 0.5462623554241067
$ws.Range("C138").Value = 0.3054124294241067
$ws.Range("C139").Value = 0.3912781354241067
$ws.Range("B140").Value = 0.0999529544241067

# --- New rows 141-145 (new ifoCAST-sampling diff dates appended) ---
# Copy formatting (bold/bordered style used by column A labels) from the
# last existing label row down into the newly appended rows.
$ws.Range("A140").Copy($ws.Range("A141:A145"))

$ws.Range("A141").Value = "2025-07-25_diff"
$ws.Range("B141").Value = 0.05603945542410671

$ws.Range("A142").Value = "2025-08-07_diff"
$ws.Range("A143").Value = "2025-08-22_diff"
$ws.Range("A144").Value = "2025-08-25_diff"
$ws.Range("A145").Value = "2025-09-08_diff"
